# Mise à jour de l'application
# Append 16 new "Wellness" entries (rows 858-873, date 2026-02-18 / serial 46071)
# to the bottom of the existing log on Feuil1, matching the format of the
# rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0), used as the trailing-space character
# throughout this workbook's "Localisation douleur" strings.
$nbsp = [char]0x00A0

# --- New row data --------------------------------------------------------
# rowNum, player name, Volume(C), Intensite(D), Fatigue(E), Douleur(F),
# Localisation douleur(G, $null = blank), Plaisir(H), expected Charge(I)
$rows = @(
    @(858, "Yoann Martelat",   70, 5, 5, 5, "Genou",               5),
    @(859, "Kamal Bafounta",   70, 5, 10, 2, "Genou psoas$nbsp",   5),
    @(860, "Maé Clavel",       70, 6, 4, 6, "Tibia",                5),
    @(861, "Jeremie Laurent",  70, 7, 7, 0, $null,                  7),
    @(862, "Naim Ighbane",     70, 3, 3, 6, "Genou",                5),
    @(863, "Omar Benyounes",   70, 3, 5, 0, $null,                  8),
    @(864, "Mehdi Boussaid",   70, 5, 5, 1, "Adducteur$nbsp",       5),
    @(865, "Levy Ndoutoume",   70, 6, 7, 0, $null,                  5),
    @(866, "Mattheo Haon",     70, 7, 6, 0, $null,                  5),
    @(867, "Emmanuel Valey",   70, 7, 8, 0, $null,                 10),
    @(868, "Ilan Ihaddadene",  70, 6, 6, 0, $null,                  7),
    @(869, "Naim Dhib",        70, 5, 5, 2, "Psoas$nbsp",           4),
    @(870, "Karahali Souaré",  70, 5, 6, 6, "Cheville",            10),
    @(871, "Theo Owono",       70, 4, 3, 3, "Coup",                 8),
    @(872, "Sofiane Belle",    70, 6, 6, 3, "Ventre",               4),
    @(873, "Romain Thunet",    70, 4, 3, 0, $null,                  1)
)

$dateSerial = 46071

# --- Formatting templates --------------------------------------------------
# Column A (date), B-F and H (name font) formats already live on row 857.
# Column G has two flavours already on the sheet -
#   G836 = empty / "Helvetica" style, G838 = filled / "Helvetica Neue" style.
$gBlankTemplate = $ws.Range("G836")
$gFilledTemplate = $ws.Range("G838")

foreach ($row in $rows) {
    $r        = $row[0]
    $name     = $row[1]
    $volume   = $row[2]
    $intens   = $row[3]
    $fatigue  = $row[4]
    $douleur  = $row[5]
    $loc      = $row[6]
    $plaisir  = $row[7]

    # Formats first (so subsequent Value assignment doesn't get clobbered).
    $ws.Range("A857").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("B857:F857").Copy() | Out-Null
    $ws.Range("B$r`:F$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("H857").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null

    if ($null -eq $loc) {
        $gBlankTemplate.Copy() | Out-Null
    } else {
        $gFilledTemplate.Copy() | Out-Null
    }
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null

    # Now the values.
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $volume
    $ws.Cells.Item($r, 4).Value = $intens
    $ws.Cells.Item($r, 5).Value = $fatigue
    $ws.Cells.Item($r, 6).Value = $douleur
    if ($null -eq $loc) {
        $ws.Cells.Item($r, 7).Value = ""
    } else {
        $ws.Cells.Item($r, 7).Value = $loc
    }
    $ws.Cells.Item($r, 8).Value = $plaisir
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"
}

$excel.CutCopyMode = $false

# Match the author's final selection / cursor position.
$ws.Range("I877").Select() | Out-Null
